$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.09982999999999999
$ws.Range("H2").Value = 0.29949
$ws.Range("I2").Value = 0.06953924013247029
$ws.Range("J2").Value = 0.06953924013247029
$ws.Range("M2").Value = 0.789222
$ws.Range("N2").Value = 2.367666
$ws.Range("O2").Value = 0.01341929863527565
$ws.Range("P2").Value = 0.01341929863527565
$ws.Range("Q2").Value = 0.07878803225999999
$ws.Range("R2").Value = 0.7090922903399999
$ws.Range("S2").Value = 0.000933167830207764
$ws.Range("T2").Value = 0.000933167830207764
$ws.Range("G3").Value = 0.09982999999999999
$ws.Range("H3").Value = 0.29949
$ws.Range("I3").Value = 0.06953924013247029
$ws.Range("J3").Value = 0.06953924013247029
$ws.Range("O3").Value = 0.005047365584441773
$ws.Range("P3").Value = 0.005047365584441773
$ws.Range("Q3").Value = 0.02963433584
$ws.Range("R3").Value = 0.26670902256
$ws.Range("S3").Value = 0.0003509899674128627
$ws.Range("T3").Value = 0.0003509899674128627
$ws.Range("G4").Value = 0.09982999999999999
$ws.Range("H4").Value = 0.29949
$ws.Range("I4").Value = 0.06953924013247029
$ws.Range("J4").Value = 0.06953924013247029
$ws.Range("M4").Value = 57.61405833333333
$ws.Range("N4").Value = 172.842175
$ws.Range("O4").Value = 0.9796232927683105
$ws.Range("P4").Value = 0.9796232927683105
$ws.Range("Q4").Value = 5.751611443416666
$ws.Range("R4").Value = 51.76450299075
$ws.Range("S4").Value = 0.06812225939517678
$ws.Range("T4").Value = 0.06812225939517678
$ws.Range("G5").Value = 0.09982999999999999
$ws.Range("H5").Value = 0.29949
$ws.Range("I5").Value = 0.06953924013247029
$ws.Range("J5").Value = 0.06953924013247029
$ws.Range("M5").Value = 0.1123343333333333
$ws.Range("N5").Value = 0.337003
$ws.Range("O5").Value = 0.001910043011972043
$ws.Range("P5").Value = 0.001910043011972043
$ws.Range("Q5").Value = 0.01121433649666666
$ws.Range("R5").Value = 0.10092902847
$ws.Range("S5").Value = 0.0001328229396728707
$ws.Range("T5").Value = 0.0001328229396728707
$ws.Range("I6").Value = 0.4393303855760352
$ws.Range("J6").Value = 0.4393303855760352
$ws.Range("M6").Value = 0.789222
$ws.Range("N6").Value = 2.367666
$ws.Range("O6").Value = 0.01341929863527565
$ws.Range("P6").Value = 0.01341929863527565
$ws.Range("Q6").Value = 0.497761789252
$ws.Range("R6").Value = 4.479856103268
$ws.Range("S6").Value = 0.005895505643595613
$ws.Range("T6").Value = 0.005895505643595613
$ws.Range("I7").Value = 0.4393303855760352
$ws.Range("J7").Value = 0.4393303855760352
$ws.Range("O7").Value = 0.005047365584441773
$ws.Range("P7").Value = 0.005047365584441773
$ws.Range("S7").Value = 0.002217461068356014
$ws.Range("T7").Value = 0.002217461068356014
$ws.Range("I8").Value = 0.4393303855760352
$ws.Range("J8").Value = 0.4393303855760352
$ws.Range("M8").Value = 57.61405833333333
$ws.Range("N8").Value = 172.842175
$ws.Range("O8").Value = 0.9796232927683105
$ws.Range("P8").Value = 0.9796232927683105
$ws.Range("Q8").Value = 36.33714818146112
$ws.Range("R8").Value = 327.03433363315
$ws.Range("S8").Value = 0.4303782789311671
$ws.Range("T8").Value = 0.4303782789311671
$ws.Range("I9").Value = 0.4393303855760352
$ws.Range("J9").Value = 0.4393303855760352
$ws.Range("M9").Value = 0.1123343333333333
$ws.Range("N9").Value = 0.337003
$ws.Range("O9").Value = 0.001910043011972043
$ws.Range("P9").Value = 0.001910043011972043
$ws.Range("Q9").Value = 0.07084918914377777
$ws.Range("R9").Value = 0.6376427022940001
$ws.Range("S9").Value = 0.0008391399329164892
$ws.Range("T9").Value = 0.0008391399329164893
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.705063
$ws.Range("H10").Value = 2.115189
$ws.Range("I10").Value = 0.4911303742914945
$ws.Range("J10").Value = 0.4911303742914945
$ws.Range("M10").Value = 0.789222
$ws.Range("N10").Value = 2.367666
$ws.Range("O10").Value = 0.01341929863527565
$ws.Range("P10").Value = 0.01341929863527565
$ws.Range("Q10").Value = 0.556451230986
$ws.Range("R10").Value = 5.008061078873999
$ws.Range("S10").Value = 0.006590625161472271
$ws.Range("T10").Value = 0.006590625161472271
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.705063
$ws.Range("H11").Value = 2.115189
$ws.Range("I11").Value = 0.4911303742914945
$ws.Range("J11").Value = 0.4911303742914945
$ws.Range("O11").Value = 0.005047365584441773
$ws.Range("P11").Value = 0.005047365584441773
$ws.Range("Q11").Value = 0.209296541424
$ws.Range("R11").Value = 1.883668872816
$ws.Range("S11").Value = 0.002478914548672896
$ws.Range("T11").Value = 0.002478914548672896
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.705063
$ws.Range("H12").Value = 2.115189
$ws.Range("I12").Value = 0.4911303742914945
$ws.Range("J12").Value = 0.4911303742914945
$ws.Range("M12").Value = 57.61405833333333
$ws.Range("N12").Value = 172.842175
$ws.Range("O12").Value = 0.9796232927683105
$ws.Range("P12").Value = 0.9796232927683105
$ws.Range("Q12").Value = 40.621540810675
$ws.Range("R12").Value = 365.593867296075
$ws.Range("S12").Value = 0.4811227544419667
$ws.Range("T12").Value = 0.4811227544419667
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.705063
$ws.Range("H13").Value = 2.115189
$ws.Range("I13").Value = 0.4911303742914945
$ws.Range("J13").Value = 0.4911303742914945
$ws.Range("M13").Value = 0.1123343333333333
$ws.Range("N13").Value = 0.337003
$ws.Range("O13").Value = 0.001910043011972043
$ws.Range("P13").Value = 0.001910043011972043
$ws.Range("Q13").Value = 0.07920278206299999
$ws.Range("R13").Value = 0.712825038567
$ws.Range("S13").Value = 0.000938080139382683
$ws.Range("T13").Value = 0.0009380801393826831
